# Apply the ErrorDocumentation updates:
# - Add a new shared string for the DataIntegrityViolationException error
# - Fill in H12 (was missing) with "n"
# - Add new rows 13 and 14 with H/I values
# - Update selection to B12 and clear the frozen/topLeftCell scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: H12 = "n" (I12 already has the Sequence "HIBERNATE_SEQUENCE" text)
$ws.Range("H12").Value = "n"

# Row 13: H13 = "y", I13 = new DataIntegrityViolationException error text
$ws.Range("H13").Value = "y"
$ws.Range("I13").Value = 'org.springframework.dao.DataIntegrityViolationException: could not execute statement; SQL [n/a]; constraint ["FK_USERROLE_USER_ID: PUBLIC.USERROLE FOREIGN KEY(USER_ID) REFERENCES PUBLIC.USER(USER_ID) (0)"; SQL statement:'

# Row 14: H14 = "n", I14 = same Sequence "HIBERNATE_SEQUENCE" text as I12
$ws.Range("H14").Value = "n"
$ws.Range("I14").Value = 'Sequence "HIBERNATE_SEQUENCE" not found; SQL statement:'

# Update the view: select B12 and scroll so topLeftCell reverts to default (A1)
$ws.Range("B12").Select() | Out-Null
